$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 56

# Column A holds a date-looking string ("2025-10-03"). Excel would normally
# auto-convert that to a real date serial number on assignment, but the
# source data stores it as plain text, so force text interpretation by
# temporarily marking the cell as Text before assigning the value, then
# restore the default (un-styled) cell appearance so no extra formatting
# is left behind on the new row - matching the rest of the sheet.
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "2025-10-03"
$ws.Cells.Item($newRow, 1).Style = "Normal"

$ws.Cells.Item($newRow, 2).Value = "21:20:50"
$ws.Cells.Item($newRow, 3).Value = "1.00 EUR = 1,793.5664"
